# Actualización automática 2025-10-14 09:30:10
# Insert a new advisor row ("ROSHANN") at row 41, alphabetically between
# "RAMIREZ MOREIRA MAYRA JACQUELINE" and "SALAZAR VERA ENRIQUE WILLIAM",
# on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. This pushes
# all following rows (and the trailing totals/summary row) down by one.

$wb = $excel.ActiveWorkbook

$sheetNames = @("VENTAS POR GRUPO", "VENTA MENSUAL")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastCol = $ws.UsedRange.Columns.Count

    # Insert a brand-new row at position 41, shifting rows 41..50 down to 42..51
    $ws.Rows.Item(41).Insert()

    # Populate the new row with the office name, the new advisor, and zeros
    $ws.Cells.Item(41, 1).Value = "OFICINA-CATAECSA"
    $ws.Cells.Item(41, 2).Value = "ROSHANN"

    for ($c = 3; $c -le $lastCol; $c++) {
        $ws.Cells.Item(41, $c).Value = 0
    }
}

# --- Sheet 1 "VENTAS POR GRUPO": the trailing "x de 48" counters now read "x de 49" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$summaryRow = 51
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item($summaryRow, $c)
    $cell.Value = $cell.Value().Replace("de 48", "de 49")
}
